$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclaimer date from 2021-07-07 to 2021-07-08
$newline = [char]10
$ws.Range("A33").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + $newline + "Model holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for holdings rows 2-29
$ws.Range("D2").Value = 0.01968964436431957
$ws.Range("E2").Value = -0.005088846250104306
$ws.Range("D3").Value = 0.01809141094757785
$ws.Range("E3").Value = 0.0002971424798192679
$ws.Range("D4").Value = 0.07554329558030148
$ws.Range("E4").Value = -0.01130667172699529
$ws.Range("D5").Value = 0.05519945514023255
$ws.Range("E5").Value = 0.009422222703147298
$ws.Range("D6").Value = 0.07124047432425994
$ws.Range("E6").Value = -0.00919969564916634
$ws.Range("D7").Value = 0.02030397046695205
$ws.Range("E7").Value = -0.02592464569650876
$ws.Range("D8").Value = 0.03383347999609961
$ws.Range("E8").Value = -0.02440251572327046
$ws.Range("D9").Value = 0.02863112047018495
$ws.Range("E9").Value = 0
$ws.Range("D10").Value = 0.02352298548747911
$ws.Range("E10").Value = -0.003497790868924877
$ws.Range("D11").Value = 0.02605883523830671
$ws.Range("E11").Value = 0.001031459515213884
$ws.Range("D12").Value = 0.02616858965868466
$ws.Range("E12").Value = -0.01380923849467885
$ws.Range("D13").Value = 0.04347350192309157
$ws.Range("E13").Value = -0.01527142327191788
$ws.Range("D14").Value = 0.02339755186419003
$ws.Range("E14").Value = -0.01018584703359549
$ws.Range("D15").Value = 0.04047563832648252
$ws.Range("E15").Value = -0.001947936957676588
$ws.Range("D16").Value = 0.0298154229300727
$ws.Range("E16").Value = -0.01725372745621467
$ws.Range("D17").Value = 0.04466870516214611
$ws.Range("E17").Value = -0.009975395806589682
$ws.Range("D18").Value = 0.1170421138910421
$ws.Range("E18").Value = -0.008966527346122244
$ws.Range("D19").Value = 0.02869921300854188
$ws.Range("E19").Value = 0.002934565434565384
$ws.Range("D20").Value = 0.02433621347846999
$ws.Range("E20").Value = -0.02302207714113902
$ws.Range("D21").Value = 0.02442536094645045
$ws.Range("E21").Value = -0.005807875479149738
$ws.Range("D22").Value = 0.01331074731710172
$ws.Range("E22").Value = -0.007000302897721467
$ws.Range("D23").Value = 0.01468991987031358
$ws.Range("E23").Value = -0.002541296060991072
$ws.Range("D24").Value = 0.03068643998322176
$ws.Range("E24").Value = -0.000145985401459936
$ws.Range("D25").Value = 0.0111295462035498
$ws.Range("E25").Value = -0.01062630816293653
$ws.Range("D26").Value = 0.03688286097370239
$ws.Range("E26").Value = -0.004153913423699307
$ws.Range("D27").Value = 0.02363587574843929
$ws.Range("E27").Value = -0.01503626393065638
$ws.Range("D28").Value = 0.05375726712389213
$ws.Range("E28").Value = -0.01412499999999994
$ws.Range("D29").Value = 0.0412903595748935
$ws.Range("E29").Value = -0.0001157273463718855

# Update Total row (row 30) Percent Change value
$ws.Range("E30").Value = -0.007987433940158439

# Restore sheet protection (it was protected before this edit, and the
# underlying diff does not alter protection settings)
$ws.Protect()
